# This script applies a cyclic "shift down by one row" edit to rows 6-11
# of the active worksheet: the data that was in row 11 moves to row 6, and
# the data that was in rows 6..10 each moves down into the next row (7..11).
#
# Only columns A, B, E, F, G, H, M, Q, R, AC actually hold differing values
# across rows 6-11 in this sheet; every other populated column (C, D, I, K,
# L, N, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AX, AY) is
# identical across all six rows, so it does not need to be rewritten.
#
# Additionally, one of the six rows (originally row 8, the "Ullticka"
# fungus record) does not have values for the K, L, M, N, AC columns at
# all (those cells are simply absent), while the other five rows (bird
# records) do have those cells (M/AC populated or blank, K/L/N blank).
# When the row data rotates, the row that ends up holding the "Ullticka"
# data must also lose those cells, while the row that used to hold it
# must gain them back (even if blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$varyCols = @("A","B","E","F","G","H","M","Q","R","AC")
$optionalCols = @("K","L","M","N","AC")
$rows = @(6,7,8,9,10,11)

# 1) Snapshot the values of the "varying" columns for every row first, so
#    that later writes do not clobber values we still need to read.
$values = @{}
foreach ($r in $rows) {
    $rowVals = @{}
    foreach ($col in $varyCols) {
        $rowVals[$col] = $ws.Range("$col$r").Value()
    }
    $values[$r] = $rowVals
}

# 2) Work out, for each destination row, which row currently supplies its
#    new content (cyclic shift down by one; row 11 wraps around to row 6).
$srcOf = @{}
foreach ($r in $rows) {
    if ($r -eq 6) {
        $srcOf[$r] = 11
    } else {
        $srcOf[$r] = $r - 1
    }
}

# 3) The sparse row (lacking K/L/M/N/AC) is currently row 8; work out where
#    its data lands after the shift.
$sparseSrcRow = 8
$sparseDestRow = 6
foreach ($r in $rows) {
    if ($srcOf[$r] -eq $sparseSrcRow) { $sparseDestRow = $r }
}

foreach ($r in $rows) {
    $src = $srcOf[$r]
    $srcVals = $values[$src]

    foreach ($col in $varyCols) {
        $v = $srcVals[$col]
        if ($null -eq $v) { $v = "" }
        $ws.Range("$col$r").Value = $v
    }

    if ($r -eq $sparseDestRow) {
        # This destination row now holds the record that has no K/L/M/N/AC
        # cells at all, so make sure those cells are removed entirely.
        foreach ($col in $optionalCols) {
            $ws.Range("$col$r").ClearContents()
        }
    } else {
        # Make sure K/L/M/N/AC cells exist (even if blank) on every other
        # destination row, matching the source row's structure. Nudging
        # the style must happen *after* the value write above, otherwise
        # writing a blank value afterwards would make the cell disappear.
        foreach ($col in $optionalCols) {
            $ws.Range("$col$r").Style = "Normal"
        }
    }
}
